{"js": "// Office.js (Word JavaScript API) edit script.\n// Goal: after the paragraph that ends with\n//   \"...restricted to CpG sites.\"\n// insert a brand-new \"BodyText\" paragraph with the DML summary sentence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\nconst marker = \"restricted to CpG sites.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the anchor paragraph (CpG sites sentence).\");\n}\n\nconst newText =\n  \"In total, 1117 differentially methylated loci (DMLs) were identified (see Methods and Supplementary Data). \" +\n  \"These DMLs were distributed across all chromosomes and classified as either hypomethylated or hypermethylated based on methylation difference.\";\n\nconst inserted = target.insertParagraph(newText, Word.InsertLocation.after);\ninserted.style = \"BodyText\";\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Goal: after the paragraph ending in \"...restricted to CpG sites.\"\n# insert a brand-new \"BodyText\" paragraph with the DML summary sentence.\n\n$d = $word.ActiveDocument\n\n$marker = \"restricted to CpG sites.\"\n$newText = \"In total, 1117 differentially methylated loci (DMLs) were identified (see Methods and Supplementary Data). These DMLs were distributed across all chromosomes and classified as either hypomethylated or hypermethylated based on methylation difference.\"\n\n# Confirm the anchor text exists before we touch the document.\n$checkRange = $d.Content\n$found = $checkRange.Find.Execute($marker)\nif (-not $found) {\n    throw \"Could not locate the anchor paragraph (CpG sites sentence).\"\n}\n\n# Locate the anchor paragraph's index within the Paragraphs collection.\n$paragraphs = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    if ($paragraphs.Item($i).Range.Text -like \"*$marker*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the anchor paragraph (CpG sites sentence).\"\n}\n\n$targetParagraph = $paragraphs.Item($targetIndex)\n$targetParagraph.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph sits right after the anchor.\n$newParagraph = $paragraphs.Item($targetIndex + 1)\n$newParagraph.Range.Text = $newText\n$newParagraph.Style = \"BodyText\"\n"}
